$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 'Archived %d note' -> 'Archived %1$d note'
$ws.Range("C14").Value = "Archived %1`$d note"

# '%d note archivée' -> '%1$d note archivée'
$ws.Range("J14").Value = "%1`$d note archivée"

# 'Archiviata %d nota' -> 'Archiviata %1$d nota'
$ws.Range("M14").Value = "Archiviata %1`$d nota"

# 'Gearchiveerde %d notitie' -> 'Gearchiveerde %1$d notitie'
$ws.Range("Q14").Value = "Gearchiveerde %1`$d notitie"

# 'Archived %d notes' -> 'Archived %1$d notes'
$ws.Range("C15").Value = "Archived %1`$d notes"

# '%d notes archivées' -> '%1$d notes archivées'
$ws.Range("J15").Value = "%1`$d notes archivées"

# 'Archiviate %d note' -> 'Archiviate %1$d note'
$ws.Range("M15").Value = "Archiviate %1`$d note"

# 'Gearchiveerde %d notities' -> 'Gearchiveerde %1$d notities'
$ws.Range("Q15").Value = "Gearchiveerde %1`$d notities"

# 'Can’t add %d file' -> 'Can’t add %1$d file'
$ws.Range("C39").Value = "Can’t add %1`$d file"

# "%d fichier n\\'a pas pu être ajouté" -> "%1$d fichier n\\'a pas pu être ajouté"
$ws.Range("J39").Value = "%1`$d fichier n\'a pas pu être ajouté"

# 'Impossibile aggiungere %d file' -> 'Impossibile aggiungere %1$d file'
$ws.Range("M39").Value = "Impossibile aggiungere %1`$d file"
$ws.Range("M40").Value = "Impossibile aggiungere %1`$d file"

# 'Kan %d bestand niet toevoegen' -> 'Kan %1$d bestand niet toevoegen'
$ws.Range("Q39").Value = "Kan %1`$d bestand niet toevoegen"

# 'Can’t add %d files' -> 'Can’t add %1$d files'
$ws.Range("C40").Value = "Can’t add %1`$d files"

# "%d fichiers n\\'ont pas pu être ajoutés" -> "%1$d fichiers n\\'ont pas pu être ajoutés"
$ws.Range("J40").Value = "%1`$d fichiers n\'ont pas pu être ajoutés"

# 'Kan %d bestanden niet toevoegen' -> 'Kan %1$d bestanden niet toevoegen'
$ws.Range("Q40").Value = "Kan %1`$d bestanden niet toevoegen"

# 'Nepodařilo se přidat %d obrázky' -> 'Nepodařilo se přidat %1$d obrázky'
$ws.Range("E43").Value = "Nepodařilo se přidat %1`$d obrázky"

# 'Nie można dodać %d obrazów' -> 'Nie można dodać %1$d obrazów'
$ws.Range("S43").Value = "Nie można dodać %1`$d obrazów"
$ws.Range("S44").Value = "Nie można dodać %1`$d obrazów"
$ws.Range("S46").Value = "Nie można dodać %1`$d obrazów"

# '%d slike niso bile dodane.' -> '%1$d slike niso bile dodane.'
$ws.Range("Y43").Value = "%1`$d slike niso bile dodane."

# 'Nepodařilo se přidat %d obrázků' -> 'Nepodařilo se přidat %1$d obrázků'
$ws.Range("E44").Value = "Nepodařilo se přidat %1`$d obrázků"

# 'Can’t add %d image' -> 'Can’t add %1$d image'
$ws.Range("C45").Value = "Can’t add %1`$d image"

# 'Nepodařilo se přidat %d obrázek' -> 'Nepodařilo se přidat %1$d obrázek'
$ws.Range("E45").Value = "Nepodařilo se přidat %1`$d obrázek"

# 'Kann %d Bild nicht hinzufügen' -> 'Kann %1$d Bild nicht hinzufügen'
$ws.Range("G45").Value = "Kann %1`$d Bild nicht hinzufügen"

# "Impossible d\\'ajouter %d image" -> "Impossible d\\'ajouter %1$d image"
$ws.Range("J45").Value = "Impossible d\'ajouter %1`$d image"

# 'Impossibile aggiungere %d immagine' -> 'Impossibile aggiungere %1$d immagine'
$ws.Range("M45").Value = "Impossibile aggiungere %1`$d immagine"

# 'Kan ikke legge til %d bilde' -> 'Kan ikke legge til %1$d bilde'
$ws.Range("P45").Value = "Kan ikke legge til %1`$d bilde"

# 'Kan %d afbeelding niet toevoegen' -> 'Kan %1$d afbeelding niet toevoegen'
$ws.Range("Q45").Value = "Kan %1`$d afbeelding niet toevoegen"

# 'Kan ikkje legga til %d bilete' -> 'Kan ikkje legga til %1$d bilete'
$ws.Range("R45").Value = "Kan ikkje legga til %1`$d bilete"
$ws.Range("R46").Value = "Kan ikkje legga til %1`$d bilete"

# 'Nie można dodać %d obrazu' -> 'Nie można dodać %1$d obrazu'
$ws.Range("S45").Value = "Nie można dodać %1`$d obrazu"

# '%d slika ni bila dodana.' -> '%1$d slika ni bila dodana.'
$ws.Range("Y45").Value = "%1`$d slika ni bila dodana."

# 'Không thể thêm %d ảnh' -> 'Không thể thêm %1$d ảnh'
$ws.Range("AD45").Value = "Không thể thêm %1`$d ảnh"
$ws.Range("AD46").Value = "Không thể thêm %1`$d ảnh"

# 'Can’t add %d images' -> 'Can’t add %1$d images'
$ws.Range("C46").Value = "Can’t add %1`$d images"

# 'Kann %d Bilder nicht hinzufügen' -> 'Kann %1$d Bilder nicht hinzufügen'
$ws.Range("G46").Value = "Kann %1`$d Bilder nicht hinzufügen"

# "Impossible d\\'ajouter %d images" -> "Impossible d\\'ajouter %1$d images"
$ws.Range("J46").Value = "Impossible d\'ajouter %1`$d images"

# 'Impossibile aggiungere %d immagini' -> 'Impossibile aggiungere %1$d immagini'
$ws.Range("M46").Value = "Impossibile aggiungere %1`$d immagini"

# 'Kan ikke legge til %d bilder' -> 'Kan ikke legge til %1$d bilder'
$ws.Range("P46").Value = "Kan ikke legge til %1`$d bilder"

# 'Kan %d afbeeldingen niet toevoegen' -> 'Kan %1$d afbeeldingen niet toevoegen'
$ws.Range("Q46").Value = "Kan %1`$d afbeeldingen niet toevoegen"

# '%d slik ni bilo dodanih.' -> '%1$d slik ni bilo dodanih.'
$ws.Range("Y46").Value = "%1`$d slik ni bilo dodanih."

# '%d sliki nista bili dodani.' -> '%1$d sliki nista bili dodani.'
$ws.Range("Y47").Value = "%1`$d sliki nista bili dodani."

# "Delete file \\'%s\\'?" -> "Delete file \\'%1$s\\'?"
$ws.Range("C84").Value = "Delete file \'%1`$s\'?"

# "Datei \\'%s\\' löschen?" -> "Datei \\'%1$s\\' löschen?"
$ws.Range("G84").Value = "Datei \'%1`$s\' löschen?"

# "Supprimer le fichier \\'%s\\'?" -> "Supprimer le fichier \\'%1$s\\'?"
$ws.Range("J84").Value = "Supprimer le fichier \'%1`$s\'?"

# 'Eliminare il file \\’%s\\’?' -> 'Eliminare il file \\’%1$s\\’?'
$ws.Range("M84").Value = "Eliminare il file \’%1`$s\’?"

# "Bestand \\'%s\\' verwijderen?" -> "Bestand \\'%1$s\\' verwijderen?"
$ws.Range("Q84").Value = "Bestand \'%1`$s\' verwijderen?"

# 'Deleted %d note' -> 'Deleted %1$d note'
$ws.Range("C93").Value = "Deleted %1`$d note"

# '%d note supprimée' -> '%1$d note supprimée'
$ws.Range("J93").Value = "%1`$d note supprimée"

# 'Eliminata %d nota' -> 'Eliminata %1$d nota'
$ws.Range("M93").Value = "Eliminata %1`$d nota"

# 'Verwijderde %d notitie' -> 'Verwijderde %1$d notitie'
$ws.Range("Q93").Value = "Verwijderde %1`$d notitie"

# 'Deleted %d notes' -> 'Deleted %1$d notes'
$ws.Range("C94").Value = "Deleted %1`$d notes"

# '%d notes supprimées' -> '%1$d notes supprimées'
$ws.Range("J94").Value = "%1`$d notes supprimées"

# 'Eliminate %d note' -> 'Eliminate %1$d note'
$ws.Range("M94").Value = "Eliminate %1`$d note"

# 'Verwijderde %d notities' -> 'Verwijderde %1$d notities'
$ws.Range("Q94").Value = "Verwijderde %1`$d notities"

# 'Imported %s Note' -> 'Imported %1$s Note'
$ws.Range("C150").Value = "Imported %1`$s Note"

# '%s note importée' -> '%1$s note importée'
$ws.Range("J150").Value = "%1`$s note importée"

# 'Importata %s nota' -> 'Importata %1$s nota'
$ws.Range("M150").Value = "Importata %1`$s nota"

# 'Geïmporteerde %s Notitie' -> 'Geïmporteerde %1$s Notitie'
$ws.Range("Q150").Value = "Geïmporteerde %1`$s Notitie"

# 'Imported %s Notes' -> 'Imported %1$s Notes'
$ws.Range("C151").Value = "Imported %1`$s Notes"

# '%s notes importées' -> '%1$s notes importées'
$ws.Range("J151").Value = "%1`$s notes importées"

# 'Importate %s note' -> 'Importate %1$s note'
$ws.Range("M151").Value = "Importate %1`$s note"

# 'Geïmporteerde %s Notities' -> 'Geïmporteerde %1$s Notities'
$ws.Range("Q151").Value = "Geïmporteerde %1`$s Notities"

# '%d more' -> '%1$d more'
$ws.Range("C190").Value = "%1`$d more"

# '%d mehr' -> '%1$d mehr'
$ws.Range("G190").Value = "%1`$d mehr"

# '%d de plus' -> '%1$d de plus'
$ws.Range("J190").Value = "%1`$d de plus"

# '…ancora %d' -> '…ancora %1$d'
$ws.Range("M190").Value = "…ancora %1`$d"

# '%d meer' -> '%1$d meer'
$ws.Range("Q190").Value = "%1`$d meer"

# '…%d more file' -> '…%1$d more file'
$ws.Range("C193").Value = "…%1`$d more file"

# '…%d weitere Datei' -> '…%1$d weitere Datei'
$ws.Range("G193").Value = "…%1`$d weitere Datei"

# '…et %d fichier ' -> '…et %1$d fichier '
$ws.Range("J193").Value = "…et %1`$d fichier "

# '…%d altro file' -> '…%1$d altro file'
$ws.Range("M193").Value = "…%1`$d altro file"

# '…%d ander bestand' -> '…%1$d ander bestand'
$ws.Range("Q193").Value = "…%1`$d ander bestand"

# '…%d more files' -> '…%1$d more files'
$ws.Range("C194").Value = "…%1`$d more files"

# '…%d weitere Dateien' -> '…%1$d weitere Dateien'
$ws.Range("G194").Value = "…%1`$d weitere Dateien"

# '…et %d fichiers' -> '…et %1$d fichiers'
$ws.Range("J194").Value = "…et %1`$d fichiers"

# '…altri %d file' -> '…altri %1$d file'
$ws.Range("M194").Value = "…altri %1`$d file"

# '…%d andere bestanden' -> '…%1$d andere bestanden'
$ws.Range("Q194").Value = "…%1`$d andere bestanden"

# 'Restored %d note' -> 'Restored %1$d note'
$ws.Range("C230").Value = "Restored %1`$d note"

# '%d note restaurée' -> '%1$d note restaurée'
$ws.Range("J230").Value = "%1`$d note restaurée"

# 'Ripristinata %d nota' -> 'Ripristinata %1$d nota'
$ws.Range("M230").Value = "Ripristinata %1`$d nota"

# 'Herstelde %d notitie' -> 'Herstelde %1$d notitie'
$ws.Range("Q230").Value = "Herstelde %1`$d notitie"

# 'Restored %d notes' -> 'Restored %1$d notes'
$ws.Range("C231").Value = "Restored %1`$d notes"

# '%d notes restaurées' -> '%1$d notes restaurées'
$ws.Range("J231").Value = "%1`$d notes restaurées"

# 'Ripristinate %d note' -> 'Ripristinate %1$d note'
$ws.Range("M231").Value = "Ripristinate %1`$d note"

# 'Herstelde %d notities' -> 'Herstelde %1$d notities'
$ws.Range("Q231").Value = "Herstelde %1`$d notities"

# 'Unarchived %d note' -> 'Unarchived %1$d note'
$ws.Range("C268").Value = "Unarchived %1`$d note"

# '%d note désarchivée' -> '%1$d note désarchivée'
$ws.Range("J268").Value = "%1`$d note désarchivée"

# 'Annullata archiviazione di %d nota' -> 'Annullata archiviazione di %1$d nota'
$ws.Range("M268").Value = "Annullata archiviazione di %1`$d nota"

# 'De-gearchiveerde %d notitie' -> 'De-gearchiveerde %1$d notitie'
$ws.Range("Q268").Value = "De-gearchiveerde %1`$d notitie"

# 'Unarchived %d notes' -> 'Unarchived %1$d notes'
$ws.Range("C269").Value = "Unarchived %1`$d notes"

# '%d notes désarchivées' -> '%1$d notes désarchivées'
$ws.Range("J269").Value = "%1`$d notes désarchivées"

# 'Annullata archiviazione di %d note' -> 'Annullata archiviazione di %1$d note'
$ws.Range("M269").Value = "Annullata archiviazione di %1`$d note"

# 'De-gearchiveerde %d notities' -> 'De-gearchiveerde %1$d notities'
$ws.Range("Q269").Value = "De-gearchiveerde %1`$d notities"
